# The target deck is a freshly (re)generated "new presentation.pptx"
# template: it ships with zero slides (the single default slide that
# PowerPoint seeds a brand new deck with is removed again).
$p = $ppt.ActivePresentation

while ($p.Slides.Count -gt 0) {
    $p.Slides.Item(1).Delete()
}
